$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column K, matching style of existing header row (A1:J1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Fill in intervention_type values for each data row (2-27)
$values = @{
    2  = "BEHAVIORAL"
    3  = "PROCEDURE"
    4  = "BEHAVIORAL"
    5  = "PROCEDURE"
    6  = "DEVICE"
    7  = "PROCEDURE"
    8  = "DEVICE"
    9  = "DEVICE"
    10 = "BEHAVIORAL"
    11 = "OTHER"
    12 = "BEHAVIORAL"
    13 = "BEHAVIORAL"
    14 = "BIOLOGICAL"
    15 = "BEHAVIORAL"
    16 = "BEHAVIORAL"
    17 = "GENETIC"
    18 = "DEVICE"
    19 = "OTHER"
    20 = "OTHER"
    21 = "BEHAVIORAL"
    22 = "BEHAVIORAL"
    23 = "BEHAVIORAL"
    24 = "OTHER"
    25 = "DIAGNOSTIC_TEST"
    26 = "OTHER"
    27 = "DEVICE"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}
